# Applies the "Merging Excel sheets is now possible" OpenTBS demo edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Accent2 (theme color index "6" in the XlThemeColor/COM enum, which maps to
# the 0-based theme="5" in the OOXML styles part), shade "darker 25%"
# (tint -0.249977111117893) -> resolved to RGB 90,3C,3A because this
# interop always serialises font colors as explicit RGB rather than as a
# theme+tint pair.
# COM long colors are BGR-ordered: B*65536 + G*256 + R
$accent2Dark25 = (0x3A * 65536) + (0x3C * 256) + 0x90

# Light gray fill ("White, Background 1, Darker 15%" -> D9D9D9)
$headerFill = (0xD9 * 65536) + (0xD9 * 256) + 0xD9

# --- Bullet-list paragraphs (B12:B16) -------------------------------------
# B12 keeps its existing text, only the colour/font changes; B13:B16 get new
# text as well.
$bullets = $ws.Range("B12:B16")
$bullets.ClearFormats()
$bullets.Font.Name = "Calibri"
$bullets.Font.Size = 11
$bullets.Font.Bold = $true
$bullets.Font.Color = $accent2Dark25

$ws.Range("B13").Value = "Merging Microsoft Excel templates with OpenTBS has several limitations because of the OpenXML format for Excel."
$ws.Range("B14").Value = "* Formulas won't work because OpenTBS needs to convert cell positions from aboslute to relative in order to have a constistent merged sheet."
$ws.Range("B15").Value = "* Formulas may also make troubles because they are saved twice in the sheet:  one for the expression, and one for the instant result."
$ws.Range("B16").Value = "* Changing picture (using ope=changepic)  because drawing information are saved in another XML sub-file."

# --- Remove the old row 17 paragraph entirely -----------------------------
$ws.Rows.Item(17).Delete()
$ws.Rows.Item(17).Insert()
$ws.Range("B17").Clear()

# --- New "Example #1" section title (B18) ---------------------------------
$ws.Range("B18").ClearFormats()
$ws.Range("B18").Value = "Example #1: merging data with rows"
$ws.Range("B18").Font.Name = "Calibri"
$ws.Range("B18").Font.Size = 11
$ws.Range("B18").Font.Bold = $true

# --- New table header row (row 20) ----------------------------------------
$header = $ws.Range("B20:D20")
$header.ClearFormats()
$ws.Range("B20").Value = "First Name"
$ws.Range("C20").Value = "Name"
$ws.Range("D20").Value = "Membership number"
$header.Interior.Color = $headerFill
$header.Borders.LineStyle = 1
$header.Borders.Weight = 2

# --- New table data row (row 21) -------------------------------------------
$dataRow = $ws.Range("B21:D21")
$dataRow.ClearFormats()
$ws.Range("B21").Value = "[a.firstname;block=row]"
$ws.Range("C21").Value = "[a.name]"
$ws.Range("D21").Value = "[a.number]"
$dataRow.Borders.LineStyle = 1
$dataRow.Borders.Weight = 2

# --- Column widths for the new table ---------------------------------------
$ws.Columns.Item(2).ColumnWidth = 15.8
$ws.Columns.Item(3).ColumnWidth = 12.1
$ws.Columns.Item(4).ColumnWidth = 18.95

# --- Selection shown when the sheet is reopened -----------------------------
[void]$ws.Range("B16").Select()

Write-Output "OpenTBS demo sheet updated"
